$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing shared string text for D3 (loyalty kickoff note)
$ws.Range("D3").Value = "Ponta-pé inicial da parte de fidelidade do sistema(Tela e definições básicas)"

# Fill in C3 with a time value (0 = midnight, matches time number format)
$ws.Range("C3").Value = 0
$ws.Range("C3").NumberFormat = "h:mm"

# Row 4: new entry
$ws.Range("A4").Value = 42622
$ws.Range("A4").NumberFormat = "d-mmm"

$ws.Range("B4").Value = 0.91666666666666663
$ws.Range("B4").NumberFormat = "h:mm"

$ws.Range("D4").Value = "Definição final sobre o funcionamento da fidelidade"

# Update active selection to D5
$ws.Range("D5").Select()
